$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Treatment query (row 5 / cell B5): remove the redundant CONCAT() wrapper
# around REPLACE(trt.treatment_agent, ';', ', ') in the "Treatment Agent" column expression.
$b5Text = @'
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs000720' AND prt.sex_at_birth = 'Male'
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
'@
$ws.Range("B5").Value = $b5Text

# Re-apply the font size on B5 (Excel records this as a fresh style entry when the
# formatting is (re)applied, matching the author's edit which bumped the cellXfs/fonts tables).
$ws.Range("B5").Font.Size = 12

# --- Adjust the saved view state: the sheet now scrolls to show row 5 at the top and
# the active selection is C5 (previously topLeftCell was A6 / selection C7).
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$ws.Range("C5").Select() | Out-Null
